$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 input values (A2:C2) - new polling figures
$ws.Range("A2").Value = 32.200000000000003
$ws.Range("B2").Value = 35.299999999999997
$ws.Range("C2").Value = 10.199999999999999

# D2 and F2 now unavailable (#N/A), matching E2/G2 which are already #N/A
$ws.Range("D2").Value = "#N/A"
$ws.Range("F2").Value = "#N/A"

# H2 becomes a formula that derives the remainder from the other parties
$ws.Range("H2").Formula = "=100-11.9-A2-B2-C2"

# Update the selected cell to reflect where the author left off
$ws.Range("J2").Select()

$wb.Save()
